$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, copying the style of O1 (bold, centered, bordered)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update swapped values (I,K,M,O) and add new columns (P,Q) for data rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O column: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P column: new, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q column: new, value 2
}
